# Applies updated crypto price/volume values to the worksheet.
# Values are written as text strings to preserve formats such as
# "309.19", "2.11%", "0.00000000751", etc. exactly as stored in the
# original inline-string cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2"  = "309.19"
    "E2"  = "2.11%"
    "D3"  = "38.90"
    "E3"  = "8.95%"
    "D4"  = "5.086"
    "E4"  = "1.06%"
    "D5"  = "0.08191"
    "E5"  = "3.62%"
    "D6"  = "2.013"
    "E6"  = "9.15%"
    "D7"  = "4.178"
    "E7"  = "1.84%"
    "D8"  = "7.913"
    "E8"  = "1.73%"
    "D9"  = "0.9342"
    "E9"  = "1.77%"
    "D10" = "0.1410"
    "E10" = "4.26%"
    "D11" = "0.1953"
    "E11" = "3.34%"
    "D12" = "0.09281"
    "E12" = "2.79%"
    "D13" = "0.03475"
    "E13" = "0.02%"
    "D14" = "0.09844"
    "E14" = "0.40%"
    "D15" = "0.001402"
    "E15" = "0.07%"
    "D16" = "0.005853"
    "E16" = "-4.31%"
    "D17" = "3.681"
    "E17" = "-1.27%"
    "E18" = "4.49%"
    "E19" = "0.43%"
    "D20" = "0.1303"
    "E20" = "-0.48%"
    "D21" = "4.811"
    "E21" = "-7.24%"
    "D22" = "0.2453"
    "E22" = "11.90%"
    "D23" = "0.04472"
    "E23" = "1.53%"
    "D24" = "0.001239"
    "E24" = "0.27%"
    "E25" = "-9.33%"
    "E27" = "-0.04%"
    "D39" = "0.02138"
    "E39" = "10.77%"
    "D40" = "0.05180"
    "E40" = "0.04%"
    "D41" = "0.007492"
    "E41" = "-1.39%"
    "D42" = "0.01013"
    "E42" = "-0.60%"
    "D43" = "0.1369"
    "E43" = "2.22%"
    "E44" = "-0.97%"
    "D45" = "0.009681"
    "E45" = "-3.15%"
    "D46" = "0.00006320"
    "E46" = "2.69%"
    "D47" = "0.00000000751"
    "E47" = "-0.03%"
    "E48" = "-0.24%"
    "E49" = "-3.56%"
    "D50" = "0.00002102"
    "E50" = "-0.03%"
    "D51" = "0.0002002"
    "E51" = "-0.03%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
